$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the Job column letters (A,B,C,D) with plain numbers (1,2,3,4)
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

# D1 header was referencing a shared string "machine3" - keep same text
$ws.Range("D1").Value = "machine3"

# Update the active selection to A5
$ws.Range("A5").Select()
